$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the "task" column first (matches the order the new shared
# strings were interned in the target file).
$ws.Range("B30").Value = "get the value from the database using java"
$ws.Range("B31").Value = "Prepar for second review "
$ws.Range("B32").Value = "edit the document for the correction."
$ws.Range("B33").Value = "insert the data from the ui page"
$ws.Range("B34").Value = "split the question from the database"
$ws.Range("B35").Value = "insert the student details from the databse"
$ws.Range("B36").Value = "edit the ui design for the test level's"
$ws.Range("B37").Value = "The css file and javascript file insert"

# Row 32's date was typed as plain text (leading spaces, dd-mm-yyyy style)
# rather than a real date value/format.
$ws.Range("A32").Value = "  10-02-2020"

# The other new dates keep the same date format as the rows above them,
# so copy that formatting down instead of re-deriving a new number format.
$ws.Range("A29").Copy()
$ws.Range("A30:A31").PasteSpecial(-4122)
$ws.Range("A33:A37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A30").Value = 43867
$ws.Range("A31").Value = 43868
$ws.Range("A33").Value = 43872
$ws.Range("A34").Value = 43873
$ws.Range("A35").Value = 43874
$ws.Range("A36").Value = 43875
$ws.Range("A37").Value = 43876

$ws.Range("A27").Select()
$excel.ActiveWindow.ScrollRow = 16
